# Generate Report for Handoff
# Marks b.md as "Ready for handoff" across the Overview / zh-cn / de-de
# sheets, records the new handoff xliff files + timestamps, and records
# the stale-handback error detail message for the new handoff round.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a72aba820306c679db589ea57ea79ae74a61647/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/142f203b209c6fe33374a1a3663a2de28467b4bc/e2e/b.md."

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-29 00:37:27"

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-29 00:37:21"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").EntireColumn.ColumnWidth = 39.14

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-29 00:37:27"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").EntireColumn.ColumnWidth = 39.14
